$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose text values look like plain numbers need to be forced to
# Text format first, otherwise Excel auto-converts them to numeric values
# (losing trailing zeros / exact decimal representation), which would not
# match the original inline-string cell type in the workbook.
function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

$ws.Range("D2").Value = "68.639.49"
$ws.Range("E2").Value = "  +2.60%  "
$ws.Range("D3").Value = "2.652.46"
$ws.Range("E3").Value = "  +2.12%  "
$ws.Range("E4").Value = "  +0.05%  "
Set-TextValue $ws.Range("D5") "599.57"
$ws.Range("E5").Value = "  +1.72%  "
Set-TextValue $ws.Range("D6") "154.73"
$ws.Range("E6").Value = "  +3.54%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("E8").Value = "  +0.36%  "
$ws.Range("D9").Value = "2.651.34"
$ws.Range("E9").Value = "  +2.23%  "
$ws.Range("E10").Value = "  +11.41%  "
$ws.Range("E11").Value = "  -0.45%  "
Set-TextValue $ws.Range("D12") "5.22"
$ws.Range("E12").Value = "  +1.31%  "
$ws.Range("E13").Value = "  +1.39%  "
$ws.Range("E14").Value = "  +2.86%  "
$ws.Range("E15").Value = "  +5.93%  "
$ws.Range("D16").Value = "3.134.35"
$ws.Range("E16").Value = "  +2.17%  "
$ws.Range("D17").Value = "68.573.07"
$ws.Range("E17").Value = "  +2.67%  "
$ws.Range("D18").Value = "2.656.26"
$ws.Range("E18").Value = "  +2.37%  "
$ws.Range("E19").Value = "  +3.90%  "
Set-TextValue $ws.Range("D20") "367.01"
$ws.Range("E20").Value = "  +1.37%  "
$ws.Range("E21").Value = "  +1.79%  "
$ws.Range("E22").Value = "  -0.14%  "
$ws.Range("E23").Value = "  +0.93%  "
$ws.Range("E24").Value = "  +4.71%  "
Set-TextValue $ws.Range("D25") "72.64"
$ws.Range("E25").Value = "  +0.57%  "
$ws.Range("E26").Value = "  +0.05%  "
$ws.Range("E27").Value = "  +0.77%  "
Set-TextValue $ws.Range("D28") "0.0000106"
$ws.Range("E28").Value = "  +7.74%  "
$ws.Range("D29").Value = "2.787.96"
$ws.Range("E29").Value = "  +2.23%  "
Set-TextValue $ws.Range("D30") "0.999"
$ws.Range("E30").Value = "  -0.13%  "
Set-TextValue $ws.Range("D31") "574.48"
$ws.Range("E31").Value = "  -0.82%  "
$ws.Range("E32").Value = "  +4.54%  "
$ws.Range("E33").Value = "  +4.81%  "
$ws.Range("E34").Value = "  +2.95%  "
$ws.Range("E35").Value = "  +4.75%  "
$ws.Range("E36").Value = "  +0.10%  "
$ws.Range("E37").Value = "  +3.89%  "
Set-TextValue $ws.Range("D38") "158.77"
$ws.Range("E38").Value = "  +1.70%  "
$ws.Range("E39").Value = "  +5.32%  "
Set-TextValue $ws.Range("D40") "19.26"
$ws.Range("E40").Value = "  +1.77%  "
$ws.Range("E41").Value = "  +3.89%  "
$ws.Range("E42").Value = "  +0.80%  "
$ws.Range("E43").Value = "  +6.46%  "
Set-TextValue $ws.Range("D44") "17.74"
$ws.Range("E44").Value = "  +4.62%  "
$ws.Range("E45").Value = "  +12.82%  "
$ws.Range("B46").Value = "USDe"
$ws.Range("C46").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
Set-TextValue $ws.Range("D46") "1.00"
$ws.Range("E46").Value = "  +0.07%  "
$ws.Range("B47").Value = "OKB"
$ws.Range("C47").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextValue $ws.Range("D47") "40.61"
$ws.Range("E47").Value = "  -0.36%  "
Set-TextValue $ws.Range("D48") "156.64"
$ws.Range("E48").Value = "  +3.13%  "
$ws.Range("E49").Value = "  +0.74%  "
$ws.Range("E50").Value = "  +2.61%  "
Set-TextValue $ws.Range("D51") "21.96"
$ws.Range("E51").Value = "  +3.14%  "
